$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.050.18"
$ws.Range("E2").Value = "  -5.31%  "
$ws.Range("D3").Value = "3.294.71"
$ws.Range("E3").Value = "  -6.63%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'522.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.49%  "
$ws.Range("D6").Value = "'173.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -13.99%  "
$ws.Range("D7").Value = "'0.602"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("D8").Value = "3.291.44"
$ws.Range("E8").Value = "  -6.48%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.605"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.66%  "
$ws.Range("D11").Value = "'56.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.38%  "
$ws.Range("D12").Value = "'0.133"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.52%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.10%  "
$ws.Range("D14").Value = "'9.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.59%  "
$ws.Range("D15").Value = "3.808.15"
$ws.Range("E15").Value = "  -6.93%  "
$ws.Range("D16").Value = "3.290.44"
$ws.Range("E16").Value = "  -6.74%  "
$ws.Range("D17").Value = "'0.116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.72%  "
$ws.Range("D18").Value = "63.938.06"
$ws.Range("E18").Value = "  -5.26%  "
$ws.Range("D19").Value = "'17.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.17%  "
$ws.Range("D20").Value = "'11.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.17%  "
$ws.Range("D21").Value = "'0.954"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.15%  "
$ws.Range("D22").Value = "'373.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.32%  "
$ws.Range("D23").Value = "'3.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.14%  "
$ws.Range("D24").Value = "'80.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("D25").Value = "'10.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.14%  "
$ws.Range("D26").Value = "'3.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'6.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'2.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.59%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'11.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.79%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'8.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.70%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'28.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.85%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'639.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.80%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.37%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'11.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("D35").Value = "'0.105"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.22%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'58.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.63%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'36.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.24%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.383"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.04%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0694"
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.927.05"
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.121"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.92%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.59%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.37%  "
$ws.Range("D46").Value = "'0.0396"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.84%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'2.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.125"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'2.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'135.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.05%  "
